# kpi_pepsico.xlsx: add a "kpi" column (C) computed from the existing
# "fecha"/"ventas" columns (A/B), and rename the old headers to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers: A1 "Rentabilidad"->"fecha", B1 "KPI_Pepsico"->"ventas",
# and add the new C1 "kpi" header.
$ws.Range("A1").Value = "fecha"
$ws.Range("B1").Value = "ventas"
$ws.Range("C1").Value = "kpi"

# New computed column: kpi = ventas(year) + ventas(prior year)/2,
# formatted as a plain percentage (0.00%).
$ws.Range("C2").NumberFormat = "0.00%"
$ws.Range("C2").Formula = "=B2+B3/2"

# Leave the selection on the new formula cell, matching the saved file.
$ws.Range("C2").Select()
